$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 6's formatting into row 7 (keeps the date-style numFmt on A7/G7
# instead of Excel re-deriving a brand new custom style from the raw value).
$ws.Range("A6:I6").Copy()
$ws.Range("A7:I7").PasteSpecial(-4122)   # xlPasteFormats

# Now fill in the new trade record's values/types for row 7.
$ws.Cells.Item(7, 1).Value = 42650.371041666665       # Date
$ws.Cells.Item(7, 2).Value = $true                    # Profitable
$ws.Cells.Item(7, 3).Value = 10212.25                 # Principle
$ws.Cells.Item(7, 4).Value = 10172.07                 # Start Principle
$ws.Cells.Item(7, 5).Value = 77.379997000000003       # BuyPrice
$ws.Cells.Item(7, 6).Value = 76.77                    # SellPrice
$ws.Cells.Item(7, 7).Value = $true                    # IsShortSell
$ws.Cells.Item(7, 8).Value = -0.79                    # Price Change %
$ws.Cells.Item(7, 9).Value = $false                   # Strong trade

# Column E (BuyPrice) widened to fit the new, longer value.
$ws.Columns.Item(5).ColumnWidth = 9
